$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1403.5
$ws.Range("I28").Value = 1299.2858
$ws.Range("K28").Value = 1299.2858
$ws.Range("M28").Value = -814.2858000000001
$ws.Range("H41").Value = 1099.2727
$ws.Range("I41").Value = 1650.6666
$ws.Range("J41").Value = 437.6
$ws.Range("K41").Value = 1650.6666
$ws.Range("L41").Value = 437.6
$ws.Range("M41").Value = -1210.6666
$ws.Range("N41").Value = -1317.6
$ws.Range("H107").Value = 1412
$ws.Range("I107").Value = 1413.5
$ws.Range("K107").Value = 1413.5
$ws.Range("M107").Value = 506.5
$ws.Range("H121").Value = 2056.2
$ws.Range("J121").Value = 2056.2
$ws.Range("L121").Value = 6168.599999999999
$ws.Range("N121").Value = -9662.599999999999
$ws.Range("H131").Value = 10758.714
$ws.Range("I131").Value = 1945
$ws.Range("K131").Value = 5835
$ws.Range("M131").Value = -795
$ws.Range("H132").Value = 1638.8269
$ws.Range("I132").Value = 1564.38
$ws.Range("K132").Value = 4693.14
$ws.Range("M132").Value = -2163.14
$ws.Range("H133").Value = 83693.25
$ws.Range("J133").Value = 83693.25
$ws.Range("L133").Value = 83693.25
$ws.Range("N133").Value = -93813.25
$ws.Range("H137").Value = 1857.8108
$ws.Range("I137").Value = 1860.1666
$ws.Range("K137").Value = 5580.4998
$ws.Range("M137").Value = -3030.4998

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 103803.91
$ws.Range("I74").Value = 127106.125
$ws.Range("K74").Value = 127106.125
$ws.Range("M74").Value = -126232.125
$ws.Range("H77").Value = 103803.91
$ws.Range("I77").Value = 127106.125
$ws.Range("K77").Value = 635530.625
$ws.Range("M77").Value = -631162.625
$ws.Range("H110").Value = 4653.8096
$ws.Range("I110").Value = 4496.1113
$ws.Range("K110").Value = 4496.1113
$ws.Range("M110").Value = -2451.1113
$ws.Range("H124").Value = 14250
$ws.Range("J124").Value = 14250
$ws.Range("L124").Value = 14250
$ws.Range("N124").Value = -24070

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 108990.69
$ws.Range("J99").Value = 262500
$ws.Range("L99").Value = 262500
$ws.Range("N99").Value = -265496

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 330
$ws.Range("I7").Value = 94.5
$ws.Range("J7").Value = 565.5
$ws.Range("K7").Value = 94.5
$ws.Range("L7").Value = 565.5
$ws.Range("M7").Value = 18.5
$ws.Range("N7").Value = -791.5
$ws.Range("H132").Value = 2293.543
$ws.Range("I132").Value = 2095.8823
$ws.Range("K132").Value = 6287.646900000001
$ws.Range("M132").Value = -3757.646900000001
$ws.Range("H141").Value = 237124.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 237124.25
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 237124.25
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -247484.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 41.166668
$ws.Range("I2").Value = 44.6
$ws.Range("J2").Value = 24
$ws.Range("K2").Value = 267.6
$ws.Range("L2").Value = 144
$ws.Range("M2").Value = -154.6
$ws.Range("N2").Value = -370
$ws.Range("H11").Value = 400
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H23").Value = 1052.5
$ws.Range("J23").Value = 1383.3334
$ws.Range("L23").Value = 4150.0002
$ws.Range("N23").Value = -4620.0002
$ws.Range("H26").Value = 1585.697
$ws.Range("I26").Value = 1092.5714
$ws.Range("K26").Value = 3277.7142
$ws.Range("M26").Value = -2989.7142
$ws.Range("H68").Value = 10002601
$ws.Range("I68").Value = 5251
$ws.Range("K68").Value = 15753
$ws.Range("M68").Value = -14942
$ws.Range("H71").Value = 10002601
$ws.Range("I71").Value = 5251
$ws.Range("K71").Value = 47259
$ws.Range("M71").Value = -43203
$ws.Range("H114").Value = 1228.238
$ws.Range("I114").Value = 861.6429000000001
$ws.Range("J114").Value = 1961.4286
$ws.Range("K114").Value = 2584.9287
$ws.Range("L114").Value = 5884.2858
$ws.Range("M114").Value = 669.0712999999996
$ws.Range("N114").Value = -12392.2858
$ws.Range("H122").Value = 1876.6471
$ws.Range("I122").Value = 1657.9231
$ws.Range("J122").Value = 2587.5
$ws.Range("K122").Value = 14921.3079
$ws.Range("L122").Value = 23287.5
$ws.Range("M122").Value = -12471.3079
$ws.Range("N122").Value = -28187.5
$ws.Range("H132").Value = 1579.0256
$ws.Range("I132").Value = 1251.2307
$ws.Range("J132").Value = 2234.6155
$ws.Range("K132").Value = 11261.0763
$ws.Range("L132").Value = 20111.5395
$ws.Range("M132").Value = -8731.076300000001
$ws.Range("N132").Value = -25171.5395

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 70000
$ws.Range("I69").Value = 70000
$ws.Range("K69").Value = 70000
$ws.Range("M69").Value = -69251
$ws.Range("H72").Value = 70000
$ws.Range("I72").Value = 70000
$ws.Range("K72").Value = 210000
$ws.Range("M72").Value = -206256
$ws.Range("H80").Value = 3250
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 3500
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -2502
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 3250
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 17500
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -12508
$ws.Range("N83").Value = -24984
$ws.Range("H102").Value = 41667430
$ws.Range("I102").Value = 778.6667
$ws.Range("J102").Value = 166667400
$ws.Range("K102").Value = 778.6667
$ws.Range("L102").Value = 166667400
$ws.Range("M102").Value = 843.3333
$ws.Range("N102").Value = -166670644
$ws.Range("H113").Value = 2597.3
$ws.Range("I113").Value = 2756.0952
$ws.Range("J113").Value = 2226.7778
$ws.Range("K113").Value = 2756.0952
$ws.Range("L113").Value = 2226.7778
$ws.Range("M113").Value = -586.0952000000002
$ws.Range("N113").Value = -6566.7778
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("H135").Value = 103716.336
$ws.Range("J135").Value = 103716.336
$ws.Range("L135").Value = 103716.336
$ws.Range("N135").Value = -113856.336
$ws.Range("H136").Value = 22154.322
$ws.Range("J136").Value = 22154.322
$ws.Range("L136").Value = 66462.966
$ws.Range("N136").Value = -71562.966

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 2157.0625
$ws.Range("J31").Value = 5441.3335
$ws.Range("L31").Value = 5441.3335
$ws.Range("N31").Value = -5937.3335
$ws.Range("H46").Value = 1197
$ws.Range("I46").Value = 1082.1666
$ws.Range("K46").Value = 1082.1666
$ws.Range("M46").Value = -894.1666
$ws.Range("H53").Value = 33061.5
$ws.Range("I53").Value = 33523
$ws.Range("J53").Value = 32600
$ws.Range("K53").Value = 33523
$ws.Range("L53").Value = 32600
$ws.Range("M53").Value = -33005
$ws.Range("N53").Value = -33636
$ws.Range("H54").Value = 30000
$ws.Range("J54").Value = 30000
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -31288
$ws.Range("H56").Value = 17277.166
$ws.Range("I56").Value = 8202
$ws.Range("J56").Value = 26352.334
$ws.Range("K56").Value = 8202
$ws.Range("L56").Value = 26352.334
$ws.Range("M56").Value = -7511
$ws.Range("N56").Value = -27734.334
$ws.Range("H58").Value = 7449.3335
$ws.Range("I58").Value = 7449.3335
$ws.Range("K58").Value = 7449.3335
$ws.Range("M58").Value = -7189.3335
$ws.Range("H82").Value = 7118.6665
$ws.Range("I82").Value = 8738.071
$ws.Range("J82").Value = 1450.75
$ws.Range("K82").Value = 8738.071
$ws.Range("L82").Value = 1450.75
$ws.Range("M82").Value = -8377.071
$ws.Range("N82").Value = -2172.75
$ws.Range("H85").Value = 7118.6665
$ws.Range("I85").Value = 8738.071
$ws.Range("J85").Value = 1450.75
$ws.Range("K85").Value = 8738.071
$ws.Range("L85").Value = 1450.75
$ws.Range("M85").Value = -7490.071
$ws.Range("N85").Value = -3946.75
$ws.Range("H122").Value = 458865.9
$ws.Range("I122").Value = 913360.0600000001
$ws.Range("J122").Value = 4371.727
$ws.Range("K122").Value = 2740080.18
$ws.Range("L122").Value = 13115.181
$ws.Range("M122").Value = -2737630.18
$ws.Range("N122").Value = -18015.181
$ws.Range("H132").Value = 3489.04
$ws.Range("I132").Value = 3149.3
$ws.Range("K132").Value = 9447.900000000001
$ws.Range("M132").Value = -6917.900000000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2565.3125
$ws.Range("I132").Value = 2419.6206
$ws.Range("K132").Value = 7258.861800000001
$ws.Range("M132").Value = -4728.861800000001
$ws.Range("H136").Value = 2056.2563
$ws.Range("I136").Value = 1086.4231
$ws.Range("J136").Value = 3995.923
$ws.Range("K136").Value = 3259.2693
$ws.Range("L136").Value = 11987.769
$ws.Range("M136").Value = -709.2692999999999
$ws.Range("N136").Value = -17087.769
